$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 12
$ws.Range("B12").Value = 1101
$ws.Range("C12").Value = 0.9

# Re-apply number format style to columns C:E (this causes Excel to reorder
# the custom cell styles, matching the diff's style index swap)
$ws.Range("C3:E12").NumberFormat = "0.00"
$ws.Range("F3").NumberFormat = "0.00"
$ws.Range("C2:E2").NumberFormat = "0.00"

# Update the selection shown when the workbook is next opened
$ws.Range("B13").Select()

$wb.Save()
